$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 15-36 (values shift as weekly data refreshes) ---
$ws.Range("D15").Value = 44575
$ws.Range("K15").Value = "Modesto"
$ws.Range("L15").Value = "Especial"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 21000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 21000
$ws.Range("Q15").Value = "`$/caja 18 kilos"
$ws.Range("R15").Value = "Región Metropolitana"
$ws.Range("S15").Value = 1167
$ws.Range("T15").Value = 18

$ws.Range("D16").Value = 44575
$ws.Range("K16").Value = "Modesto"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("Q16").Value = "`$/caja 18 kilos"
$ws.Range("R16").Value = "Región Metropolitana"
$ws.Range("S16").Value = 1000
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 44575
$ws.Range("K17").Value = "Modesto"
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("Q17").Value = "`$/caja 18 kilos"
$ws.Range("R17").Value = "Región Metropolitana"
$ws.Range("S17").Value = 889
$ws.Range("T17").Value = 18

$ws.Range("D18").Value = 44551
$ws.Range("K18").Value = "Castle Brite"
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 20000
$ws.Range("Q18").Value = "`$/caja 18 kilos"
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 1111
$ws.Range("T18").Value = 18

$ws.Range("D19").Value = 44551
$ws.Range("K19").Value = "Castle Brite"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("Q19").Value = "`$/caja 18 kilos"
$ws.Range("R19").Value = "Región Metropolitana"
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18

$ws.Range("D20").Value = 44551
$ws.Range("K20").Value = "Castle Brite"
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 16000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 16000
$ws.Range("Q20").Value = "`$/caja 18 kilos"
$ws.Range("R20").Value = "Región Metropolitana"
$ws.Range("S20").Value = 889
$ws.Range("T20").Value = 18

$ws.Range("D21").Value = 44537
$ws.Range("K21").Value = "Castle Brite"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 500
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 21000
$ws.Range("Q21").Value = "`$/caja 18 kilos"
$ws.Range("R21").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S21").Value = 1167
$ws.Range("T21").Value = 18

$ws.Range("D22").Value = 44537
$ws.Range("K22").Value = "Castle Brite"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 17000
$ws.Range("O22").Value = 17000
$ws.Range("P22").Value = 17000
$ws.Range("Q22").Value = "`$/caja 18 kilos"
$ws.Range("R22").Value = "Región del Maule"
$ws.Range("S22").Value = 944
$ws.Range("T22").Value = 18

$ws.Range("D23").Value = 44159
$ws.Range("K23").Value = "Castle Brite"
$ws.Range("L23").Value = "Tercera"
$ws.Range("M23").Value = 400
$ws.Range("N23").Value = 15500
$ws.Range("O23").Value = 16000
$ws.Range("P23").Value = 15750
$ws.Range("Q23").Value = "`$/caja 15 kilos"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 1050
$ws.Range("T23").Value = 15

$ws.Range("D24").Value = 44572
$ws.Range("K24").Value = "Modesto"
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 150
$ws.Range("N24").Value = 21000
$ws.Range("O24").Value = 21000
$ws.Range("P24").Value = 21000
$ws.Range("Q24").Value = "`$/caja 18 kilos"
$ws.Range("R24").Value = "Región Metropolitana"
$ws.Range("S24").Value = 1167
$ws.Range("T24").Value = 18

$ws.Range("D25").Value = 44572
$ws.Range("K25").Value = "Modesto"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 150
$ws.Range("N25").Value = 18000
$ws.Range("O25").Value = 18000
$ws.Range("P25").Value = 18000
$ws.Range("Q25").Value = "`$/caja 18 kilos"
$ws.Range("R25").Value = "Región Metropolitana"
$ws.Range("S25").Value = 1000
$ws.Range("T25").Value = 18

$ws.Range("D26").Value = 44572
$ws.Range("K26").Value = "Modesto"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 150
$ws.Range("N26").Value = 16000
$ws.Range("O26").Value = 16000
$ws.Range("P26").Value = 16000
$ws.Range("Q26").Value = "`$/caja 18 kilos"
$ws.Range("R26").Value = "Región Metropolitana"
$ws.Range("S26").Value = 889
$ws.Range("T26").Value = 18

$ws.Range("D27").Value = 44166
$ws.Range("K27").Value = "Castle Brite"
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 600
$ws.Range("N27").Value = 16000
$ws.Range("O27").Value = 17000
$ws.Range("P27").Value = 16500
$ws.Range("Q27").Value = "`$/caja 15 kilos"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 1100
$ws.Range("T27").Value = 15

$ws.Range("D28").Value = 44530
$ws.Range("K28").Value = "Castle Brite"
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 500
$ws.Range("N28").Value = 20000
$ws.Range("O28").Value = 21000
$ws.Range("P28").Value = 20500
$ws.Range("Q28").Value = "`$/caja 18 kilos"
$ws.Range("R28").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S28").Value = 1139
$ws.Range("T28").Value = 18

$ws.Range("D29").Value = 44187
$ws.Range("K29").Value = "Castle Brite"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 350
$ws.Range("N29").Value = 16000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 16000
$ws.Range("Q29").Value = "`$/caja 15 kilos"
$ws.Range("R29").Value = "Región Metropolitana"
$ws.Range("S29").Value = 1067
$ws.Range("T29").Value = 15

$ws.Range("D30").Value = 44187
$ws.Range("K30").Value = "Castle Brite"
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 13000
$ws.Range("O30").Value = 13000
$ws.Range("P30").Value = 13000
$ws.Range("Q30").Value = "`$/caja 15 kilos"
$ws.Range("R30").Value = "Región Metropolitana"
$ws.Range("S30").Value = 867
$ws.Range("T30").Value = 15

$ws.Range("D31").Value = 44540
$ws.Range("K31").Value = "Castle Brite"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 600
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("Q31").Value = "`$/caja 18 kilos"
$ws.Range("R31").Value = "Región del Maule"
$ws.Range("S31").Value = 889
$ws.Range("T31").Value = 18

$ws.Range("D32").Value = 44544
$ws.Range("K32").Value = "Castle Brite"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 600
$ws.Range("N32").Value = 18000
$ws.Range("O32").Value = 20000
$ws.Range("P32").Value = 19000
$ws.Range("Q32").Value = "`$/caja 18 kilos"
$ws.Range("R32").Value = "Región Metropolitana"
$ws.Range("S32").Value = 1056
$ws.Range("T32").Value = 18

$ws.Range("D33").Value = 44544
$ws.Range("K33").Value = "Castle Brite"
$ws.Range("L33").Value = "Segunda"
$ws.Range("M33").Value = 300
$ws.Range("N33").Value = 16000
$ws.Range("O33").Value = 16000
$ws.Range("P33").Value = 16000
$ws.Range("Q33").Value = "`$/caja 18 kilos"
$ws.Range("R33").Value = "Región Metropolitana"
$ws.Range("S33").Value = 889
$ws.Range("T33").Value = 18

$ws.Range("D34").Value = 44169
$ws.Range("K34").Value = "Castle Brite"
$ws.Range("L34").Value = "Segunda"
$ws.Range("M34").Value = 500
$ws.Range("N34").Value = 15000
$ws.Range("O34").Value = 16000
$ws.Range("P34").Value = 15500
$ws.Range("Q34").Value = "`$/caja 15 kilos"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 1033
$ws.Range("T34").Value = 15

$ws.Range("D35").Value = 44194
$ws.Range("K35").Value = "Castle Brite"
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 300
$ws.Range("N35").Value = 15000
$ws.Range("O35").Value = 16000
$ws.Range("P35").Value = 15500
$ws.Range("Q35").Value = "`$/caja 15 kilos"
$ws.Range("R35").Value = "Región Metropolitana"
$ws.Range("S35").Value = 1033
$ws.Range("T35").Value = 15

$ws.Range("D36").Value = 44162
$ws.Range("K36").Value = "Castle Brite"
$ws.Range("L36").Value = "Tercera"
$ws.Range("M36").Value = 500
$ws.Range("N36").Value = 15000
$ws.Range("O36").Value = 16000
$ws.Range("P36").Value = 15500
$ws.Range("Q36").Value = "`$/caja 15 kilos"
$ws.Range("R36").Value = "Región de O'Higgins"
$ws.Range("S36").Value = 1033
$ws.Range("T36").Value = 15

# --- Append new rows 37-39 (full rows, including the columns that stay constant) ---
$ws.Range("A37").Value = 4
$ws.Range("B37").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C37").Value = "Los Lagos"
$ws.Range("D37").Value = 44533
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100103
$ws.Range("H37").Value = "Frutos de hueso (carozo)"
$ws.Range("I37").Value = 100103003
$ws.Range("J37").Value = "Damasco"
$ws.Range("K37").Value = "Castle Brite"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 350
$ws.Range("N37").Value = 24000
$ws.Range("O37").Value = 24000
$ws.Range("P37").Value = 24000
$ws.Range("Q37").Value = "`$/caja 18 kilos"
$ws.Range("R37").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S37").Value = 1333
$ws.Range("T37").Value = 18

$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"
$ws.Range("D38").Value = 44533
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100103
$ws.Range("H38").Value = "Frutos de hueso (carozo)"
$ws.Range("I38").Value = 100103003
$ws.Range("J38").Value = "Damasco"
$ws.Range("K38").Value = "Castle Brite"
$ws.Range("L38").Value = "Segunda"
$ws.Range("M38").Value = 350
$ws.Range("N38").Value = 20000
$ws.Range("O38").Value = 20000
$ws.Range("P38").Value = 20000
$ws.Range("Q38").Value = "`$/caja 18 kilos"
$ws.Range("R38").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S38").Value = 1111
$ws.Range("T38").Value = 18

$ws.Range("A39").Value = 4
$ws.Range("B39").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C39").Value = "Los Lagos"
$ws.Range("D39").Value = 44533
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103003
$ws.Range("J39").Value = "Damasco"
$ws.Range("K39").Value = "Castle Brite"
$ws.Range("L39").Value = "Tercera"
$ws.Range("M39").Value = 350
$ws.Range("N39").Value = 17000
$ws.Range("O39").Value = 17000
$ws.Range("P39").Value = 17000
$ws.Range("Q39").Value = "`$/caja 18 kilos"
$ws.Range("R39").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S39").Value = 944
$ws.Range("T39").Value = 18

